$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# Sheet1: tweak a handful of values and append a lone C9 cell
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = 15.1
$ws1.Range("C2").Value = 2
$ws1.Range("C3").Value = 34
$ws1.Range("C4").Value = 523
$ws1.Range("B5").Value = 3
$ws1.Range("C5").Value = 45
$ws1.Range("C6").Value = 23
$ws1.Range("C7").Value = 45
$ws1.Range("C8").Value = 23
$ws1.Range("C9").Value = 54

# ---------------------------------------------------------------------
# Sheet2: tweak column C values
# ---------------------------------------------------------------------
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 2
$ws2.Range("C4").Value = 3
$ws2.Range("C5").Value = 123
$ws2.Range("C6").Value = 12
$ws2.Range("C7").Value = 3
$ws2.Range("C8").Value = 123

# ---------------------------------------------------------------------
# Sheet3: add three new columns (value / x2 / y2) plus an "asdf" column
# inserted before them, populate all of the new data, and drop the
# old, now-stale column C entries for rows 5 & 7.
# ---------------------------------------------------------------------
$ws3.Range("D1").Value = "value"
$ws3.Range("E1").Value = "x2"
$ws3.Range("F1").Value = "y2"
$ws3.Range("C1").Value = "asdf"

$ws3.Range("C2").Value = 123
$ws3.Range("D2").Value = 1
$ws3.Range("E2").Value = 1
$ws3.Range("F2").Value = 13

$ws3.Range("C3").Value = 123
$ws3.Range("D3").Value = 123
$ws3.Range("E3").Value = 2
$ws3.Range("F3").Value = 21

$ws3.Range("C4").Value = 12
$ws3.Range("D4").Value = 23
$ws3.Range("E4").Value = 3
$ws3.Range("F4").Value = 2

$ws3.Range("C5").ClearContents()
$ws3.Range("D5").Value = 4
$ws3.Range("E5").Value = 4
$ws3.Range("F5").Value = 2

$ws3.Range("C6").Value = 123412
$ws3.Range("D6").Value = 34
$ws3.Range("E6").Value = 5
$ws3.Range("F6").Value = 3

$ws3.Range("C7").ClearContents()
$ws3.Range("D7").Value = 3
$ws3.Range("E7").Value = 7
$ws3.Range("F7").Value = 3

$ws3.Range("C8").Value = 34
$ws3.Range("D8").Value = 3
$ws3.Range("E8").Value = 6
$ws3.Range("F8").Value = 4

$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Selections / active cells, applied last so the final active sheet
# matches the original (Sheet3).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D21").Select()

$ws2.Activate()
$ws2.Range("C9").Select()

$ws3.Activate()
$ws3.Range("D9").Select()
